# Auto-generated edit script: updates computed marketboard profit
# columns (H-N) for various Leve rows across multiple sheets, per
# the scheduled-runner price refresh.

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")

# Row 12
$ws.Range("H12").Value = 379.75
$ws.Range("I12").Value = 372
$ws.Range("J12").Value = 387.5
$ws.Range("K12").Value = 372
$ws.Range("L12").Value = 387.5
$ws.Range("M12").Value = -202
$ws.Range("N12").Value = -727.5
# Row 28
$ws.Range("H28").Value = 8105.8965
$ws.Range("I28").Value = 10784.571
$ws.Range("J28").Value = 1074.375
$ws.Range("K28").Value = 10784.571
$ws.Range("L28").Value = 1074.375
$ws.Range("M28").Value = -10299.571
$ws.Range("N28").Value = -2044.375
# Row 61
$ws.Range("H61").Value = 1508.5714
$ws.Range("I61").Value = 93.333336
$ws.Range("K61").Value = 280.000008
$ws.Range("M61").Value = -108.000008
# Row 107
$ws.Range("H107").Value = 277.8
$ws.Range("I107").Value = 277.8
$ws.Range("K107").Value = 277.8
$ws.Range("M107").Value = 1642.2
# Row 113
$ws.Range("H113").Value = 3016.647
$ws.Range("I113").Value = 1620
$ws.Range("J113").Value = 3778.4546
$ws.Range("K113").Value = 1620
$ws.Range("L113").Value = 3778.4546
$ws.Range("M113").Value = 1634
$ws.Range("N113").Value = -10286.4546
# Row 116
$ws.Range("H116").Value = 3911.889
$ws.Range("I116").Value = 5903
$ws.Range("J116").Value = 3146.077
$ws.Range("K116").Value = 5903
$ws.Range("L116").Value = 3146.077
$ws.Range("M116").Value = -2461
$ws.Range("N116").Value = -10030.077
# Row 132
$ws.Range("H132").Value = 3374.689
$ws.Range("I132").Value = 2965.9666
$ws.Range("J132").Value = 4192.1333
$ws.Range("K132").Value = 8897.899800000001
$ws.Range("L132").Value = 12576.3999
$ws.Range("M132").Value = -6367.899800000001
$ws.Range("N132").Value = -17636.3999
# Row 141
$ws.Range("H141").Value = 8984.565000000001
$ws.Range("I141").Value = 6076.4287
$ws.Range("J141").Value = 13508.333
$ws.Range("K141").Value = 18229.2861
$ws.Range("L141").Value = 40524.999
$ws.Range("M141").Value = -13049.2861
$ws.Range("N141").Value = -50884.999

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 359823.66
$ws.Range("I32").Value = 389484.3
$ws.Range("K32").Value = 389484.3
$ws.Range("M32").Value = -389197.3
# Row 132
$ws.Range("H132").Value = 8369.742
$ws.Range("I132").Value = 5500.4
$ws.Range("J132").Value = 20325.334
$ws.Range("K132").Value = 16501.2
$ws.Range("L132").Value = 60976.00199999999
$ws.Range("M132").Value = -13971.2
$ws.Range("N132").Value = -66036.00199999999

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")

# Row 76
$ws.Range("H76").Value = 333345440
$ws.Range("J76").Value = 333345440
$ws.Range("L76").Value = 333345440
$ws.Range("N76").Value = -333346070
# Row 79
$ws.Range("H79").Value = 333345440
$ws.Range("J79").Value = 333345440
$ws.Range("L79").Value = 333345440
$ws.Range("N79").Value = -333347624
# Row 86
$ws.Range("H86").Value = 2032
$ws.Range("I86").Value = 1580.875
$ws.Range("J86").Value = 2433
$ws.Range("K86").Value = 1580.875
$ws.Range("L86").Value = 2433
$ws.Range("M86").Value = -457.875
$ws.Range("N86").Value = -4679
# Row 89
$ws.Range("H89").Value = 2032
$ws.Range("I89").Value = 1580.875
$ws.Range("J89").Value = 2433
$ws.Range("K89").Value = 7904.375
$ws.Range("L89").Value = 12165
$ws.Range("M89").Value = -2288.375
$ws.Range("N89").Value = -23397
# Row 134
$ws.Range("H134").Value = 918.871
$ws.Range("I134").Value = 918.871
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 2756.613
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -221.6129999999998
$ws.Range("N134").ClearContents()

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")

# Row 7
$ws.Range("H7").Value = 76.71429000000001
$ws.Range("I7").Value = 56.166668
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 56.166668
$ws.Range("L7").Value = 200
$ws.Range("M7").Value = 56.833332
$ws.Range("N7").Value = -426
# Row 16
$ws.Range("H16").Value = 1853.1852
$ws.Range("I16").Value = 988.0714
$ws.Range("J16").Value = 2784.8462
$ws.Range("K16").Value = 988.0714
$ws.Range("L16").Value = 2784.8462
$ws.Range("M16").Value = -701.0714
$ws.Range("N16").Value = -3358.8462
# Row 31
$ws.Range("H31").Value = 16003.1
$ws.Range("I31").Value = 1109.6666
$ws.Range("J31").Value = 18631.354
$ws.Range("K31").Value = 1109.6666
$ws.Range("L31").Value = 18631.354
$ws.Range("M31").Value = -814.6666
$ws.Range("N31").Value = -19221.354
# Row 34
$ws.Range("H34").Value = 16003.1
$ws.Range("I34").Value = 1109.6666
$ws.Range("J34").Value = 18631.354
$ws.Range("K34").Value = 1109.6666
$ws.Range("L34").Value = 18631.354
$ws.Range("M34").Value = -907.6666
$ws.Range("N34").Value = -19035.354
# Row 74
$ws.Range("H74").Value = 18542.666
$ws.Range("J74").Value = 18542.666
$ws.Range("L74").Value = 18542.666
$ws.Range("N74").Value = -20290.666
# Row 77
$ws.Range("H77").Value = 18542.666
$ws.Range("J77").Value = 18542.666
$ws.Range("L77").Value = 55627.99800000001
$ws.Range("N77").Value = -64363.99800000001
# Row 94
$ws.Range("H94").Value = 689.38464
$ws.Range("I94").Value = 500
$ws.Range("J94").Value = 723.8182
$ws.Range("K94").Value = 500
$ws.Range("L94").Value = 723.8182
$ws.Range("M94").Value = -49
$ws.Range("N94").Value = -1625.8182
# Row 113
$ws.Range("H113").Value = 1853.1852
$ws.Range("I113").Value = 988.0714
$ws.Range("J113").Value = 2784.8462
$ws.Range("K113").Value = 988.0714
$ws.Range("L113").Value = 2784.8462
$ws.Range("M113").Value = 1181.9286
$ws.Range("N113").Value = -7124.8462
# Row 132
$ws.Range("H132").Value = 76937160
$ws.Range("I132").Value = 111125570
$ws.Range("J132").Value = 13253
$ws.Range("K132").Value = 333376710
$ws.Range("L132").Value = 39759
$ws.Range("M132").Value = -333374180
$ws.Range("N132").Value = -44819
# Row 134
$ws.Range("H134").Value = 2484.9546
$ws.Range("I134").Value = 2460.4119
$ws.Range("J134").Value = 2568.4
$ws.Range("K134").Value = 7381.2357
$ws.Range("L134").Value = 7705.200000000001
$ws.Range("M134").Value = -4846.2357
$ws.Range("N134").Value = -12775.2

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")

# Row 87
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("N87").ClearContents()
# Row 90
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("N90").ClearContents()
# Row 109
$ws.Range("H109").Value = 399.25
$ws.Range("I109").Value = 399
$ws.Range("K109").Value = 1197
$ws.Range("M109").Value = -157
# Row 113
$ws.Range("H113").Value = 353097.12
$ws.Range("I113").Value = 423.7931
$ws.Range("J113").Value = 622242.5600000001
$ws.Range("K113").Value = 1271.3793
$ws.Range("L113").Value = 1866727.68
$ws.Range("M113").Value = 898.6206999999999
$ws.Range("N113").Value = -1871067.68

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")

# Row 2
$ws.Range("H2").Value = 38.944443
$ws.Range("I2").Value = 22.636364
$ws.Range("J2").Value = 64.57143000000001
$ws.Range("K2").Value = 22.636364
$ws.Range("L2").Value = 64.57143000000001
$ws.Range("M2").Value = 90.363636
$ws.Range("N2").Value = -290.57143
# Row 70
$ws.Range("H70").Value = 8264.538
$ws.Range("I70").Value = 15794.75
$ws.Range("J70").Value = 4917.778
$ws.Range("K70").Value = 15794.75
$ws.Range("L70").Value = 4917.778
$ws.Range("M70").Value = -15524.75
$ws.Range("N70").Value = -5457.778
# Row 73
$ws.Range("H73").Value = 8264.538
$ws.Range("I73").Value = 15794.75
$ws.Range("J73").Value = 4917.778
$ws.Range("K73").Value = 15794.75
$ws.Range("L73").Value = 4917.778
$ws.Range("M73").Value = -14858.75
$ws.Range("N73").Value = -6789.778
# Row 80
$ws.Range("H80").Value = 2856.5264
$ws.Range("I80").Value = 2317.2727
$ws.Range("J80").Value = 3598
$ws.Range("K80").Value = 2317.2727
$ws.Range("L80").Value = 3598
$ws.Range("M80").Value = -1319.2727
$ws.Range("N80").Value = -5594
# Row 83
$ws.Range("H83").Value = 2856.5264
$ws.Range("I83").Value = 2317.2727
$ws.Range("J83").Value = 3598
$ws.Range("K83").Value = 11586.3635
$ws.Range("L83").Value = 17990
$ws.Range("M83").Value = -6594.363499999999
$ws.Range("N83").Value = -27974
# Row 97
$ws.Range("H97").Value = 430.14285
$ws.Range("I97").Value = 300
$ws.Range("J97").Value = 603.6667
$ws.Range("K97").Value = 300
$ws.Range("L97").Value = 603.6667
$ws.Range("M97").Value = 196
$ws.Range("N97").Value = -1595.6667
# Row 107
$ws.Range("H107").Value = 117.42857
$ws.Range("I107").Value = 117
$ws.Range("K107").Value = 117
$ws.Range("M107").Value = 1803
# Row 122
$ws.Range("H122").Value = 1126.1818
$ws.Range("I122").Value = 1041.1428
$ws.Range("J122").Value = 1275
$ws.Range("K122").Value = 3123.4284
$ws.Range("L122").Value = 3825
$ws.Range("M122").Value = -673.4284000000002
$ws.Range("N122").Value = -8725
# Row 126
$ws.Range("H126").Value = 1573.9131
$ws.Range("I126").Value = 1026.3572
$ws.Range("K126").Value = 3079.0716
$ws.Range("M126").Value = -609.0715999999998

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")

# Row 22
$ws.Range("H22").Value = 5477
$ws.Range("I22").Value = 502.92307
$ws.Range("J22").Value = 14714.571
$ws.Range("K22").Value = 502.92307
$ws.Range("L22").Value = 14714.571
$ws.Range("M22").Value = -207.92307
$ws.Range("N22").Value = -15304.571
# Row 27
$ws.Range("H27").Value = 5477
$ws.Range("I27").Value = 502.92307
$ws.Range("J27").Value = 14714.571
$ws.Range("K27").Value = 502.92307
$ws.Range("L27").Value = 14714.571
$ws.Range("M27").Value = -395.92307
$ws.Range("N27").Value = -14928.571
# Row 132
$ws.Range("H132").Value = 6869.9165
$ws.Range("I132").Value = 6415.769
$ws.Range("K132").Value = 19247.307
$ws.Range("M132").Value = -16717.307

